$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 422.8
$ws.Range("I6").Value = 278.5
$ws.Range("K6").Value = 835.5
$ws.Range("M6").Value = -723.5

$ws.Range("H8").Value = 32
$ws.Range("I8").Value = 32
$ws.Range("K8").Value = 96
$ws.Range("M8").Value = 43

$ws.Range("H15").Value = 1954.2931
$ws.Range("I15").Value = 1954.2931
$ws.Range("K15").Value = 5862.879300000001
$ws.Range("M15").Value = -5693.879300000001

$ws.Range("H17").Value = 2954.5151
$ws.Range("J17").Value = 2954.5151
$ws.Range("L17").Value = 8863.5453
$ws.Range("N17").Value = -9199.5453

$ws.Range("H21").Value = 20000
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = ""

$ws.Range("H23").Value = 20000
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = ""

$ws.Range("H34").Value = 7690
$ws.Range("I34").Value = 7690
$ws.Range("K34").Value = 7690
$ws.Range("M34").Value = -7487

$ws.Range("H36").Value = 7690
$ws.Range("I36").Value = 7690
$ws.Range("K36").Value = 7690
$ws.Range("M36").Value = -6975

$ws.Range("H80").Value = 1332.2727
$ws.Range("I80").Value = 471
$ws.Range("J80").Value = 1585.5883
$ws.Range("K80").Value = 1413
$ws.Range("L80").Value = 4756.7649
$ws.Range("M80").Value = -415
$ws.Range("N80").Value = -6752.7649

$ws.Range("H83").Value = 1332.2727
$ws.Range("I83").Value = 471
$ws.Range("J83").Value = 1585.5883
$ws.Range("K83").Value = 4239
$ws.Range("L83").Value = 14270.2947
$ws.Range("M83").Value = 753
$ws.Range("N83").Value = -24254.2947

$ws.Range("H96").Value = 1544.8636
$ws.Range("I96").Value = 433.5263
$ws.Range("K96").Value = 1300.5789
$ws.Range("M96").Value = 72.42110000000002

$ws.Range("H111").Value = 2766.8572
$ws.Range("I111").Value = 3138.1428
$ws.Range("J111").Value = 2024.2858
$ws.Range("K111").Value = 9414.428400000001
$ws.Range("L111").Value = 6072.857400000001
$ws.Range("M111").Value = -6347.428400000001
$ws.Range("N111").Value = -12206.8574

$ws.Range("H112").Value = 1702.6
$ws.Range("J112").Value = 1738.5
$ws.Range("L112").Value = 5215.5
$ws.Range("N112").Value = -7431.5

$ws.Range("H132").Value = 1379.2
$ws.Range("I132").Value = 1424
$ws.Range("K132").Value = 4272
$ws.Range("M132").Value = -1742

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 34696.973
$ws.Range("I32").Value = 35993.94
$ws.Range("K32").Value = 35993.94
$ws.Range("M32").Value = -35706.94

$ws.Range("H35").Value = 1429
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").Value = ""

$ws.Range("H61").Value = 12579.833
$ws.Range("I61").Value = 12013.357
$ws.Range("J61").Value = 14562.5
$ws.Range("K61").Value = 12013.357
$ws.Range("L61").Value = 14562.5
$ws.Range("M61").Value = -11801.357
$ws.Range("N61").Value = -14986.5

$ws.Range("H63").Value = 6446.875
$ws.Range("I63").Value = 1884.6428
$ws.Range("K63").Value = 1884.6428
$ws.Range("M63").Value = -1198.6428

$ws.Range("H66").Value = 6446.875
$ws.Range("I66").Value = 1884.6428
$ws.Range("K66").Value = 9423.214
$ws.Range("M66").Value = -5991.214

$ws.Range("H88").Value = 1337.4375
$ws.Range("I88").Value = 1334.5
$ws.Range("J88").Value = 1338.4166
$ws.Range("K88").Value = 1334.5
$ws.Range("L88").Value = 1338.4166
$ws.Range("M88").Value = -928.5
$ws.Range("N88").Value = -2150.4166

$ws.Range("H91").Value = 1337.4375
$ws.Range("I91").Value = 1334.5
$ws.Range("J91").Value = 1338.4166
$ws.Range("K91").Value = 1334.5
$ws.Range("L91").Value = 1338.4166
$ws.Range("M91").Value = 69.5
$ws.Range("N91").Value = -4146.4166

$ws.Range("H119").Value = 50000
$ws.Range("J119").Value = 50000
$ws.Range("L119").Value = 50000
$ws.Range("N119").Value = -59676

$ws.Range("H122").Value = 3814.2144
$ws.Range("I122").Value = 3360
$ws.Range("K122").Value = 10080
$ws.Range("M122").Value = -7630

$ws.Range("H128").Value = 198823
$ws.Range("J128").Value = 198823
$ws.Range("L128").Value = 198823
$ws.Range("N128").Value = -208783

$ws.Range("H130").Value = 159290.6
$ws.Range("J130").Value = 159290.6
$ws.Range("L130").Value = 159290.6
$ws.Range("N130").Value = -169330.6

$ws.Range("H132").Value = 4053.814
$ws.Range("I132").Value = 3125
$ws.Range("J132").Value = 7562.6665
$ws.Range("K132").Value = 9375
$ws.Range("L132").Value = 22687.9995
$ws.Range("M132").Value = -6845
$ws.Range("N132").Value = -27747.9995

$ws.Range("H136").Value = 12579.833
$ws.Range("I136").Value = 12013.357
$ws.Range("J136").Value = 14562.5
$ws.Range("K136").Value = 36040.071
$ws.Range("L136").Value = 43687.5
$ws.Range("M136").Value = -33490.071
$ws.Range("N136").Value = -48787.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 33163.06
$ws.Range("J88").Value = 33163.06
$ws.Range("L88").Value = 33163.06
$ws.Range("N88").Value = -33975.06

$ws.Range("H91").Value = 33163.06
$ws.Range("J91").Value = 33163.06
$ws.Range("L91").Value = 33163.06
$ws.Range("N91").Value = -35971.06

$ws.Range("H107").Value = 1150.24
$ws.Range("I107").Value = 1063.7333
$ws.Range("J107").Value = 1280
$ws.Range("K107").Value = 1063.7333
$ws.Range("L107").Value = 1280
$ws.Range("M107").Value = 856.2666999999999
$ws.Range("N107").Value = -5120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 50042.46
$ws.Range("I20").Value = 82000.336
$ws.Range("J20").Value = 40455.1
$ws.Range("K20").Value = 82000.336
$ws.Range("L20").Value = 40455.1
$ws.Range("M20").Value = -81755.336
$ws.Range("N20").Value = -40945.1

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = ""

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = ""

$ws.Range("H132").Value = 6918
$ws.Range("J132").Value = 12499.5
$ws.Range("L132").Value = 37498.5
$ws.Range("N132").Value = -42558.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1144.2858
$ws.Range("I9").Value = 219.8
$ws.Range("J9").Value = 3455.5
$ws.Range("K9").Value = 219.8
$ws.Range("L9").Value = 3455.5
$ws.Range("M9").Value = 4.199999999999989
$ws.Range("N9").Value = -3903.5

$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").Value = ""

$ws.Range("H35").Value = 2311.75
$ws.Range("I35").Value = 2311.75
$ws.Range("K35").Value = 2311.75
$ws.Range("M35").Value = -1975.75

$ws.Range("H40").Value = 8478.947
$ws.Range("I40").Value = 5706.933
$ws.Range("K40").Value = 5706.933
$ws.Range("M40").Value = -5570.933

$ws.Range("H61").Value = 1736.2142
$ws.Range("I61").Value = 1370.3
$ws.Range("K61").Value = 1370.3
$ws.Range("M61").Value = -1168.3

$ws.Range("H113").Value = 1736.2142
$ws.Range("I113").Value = 1370.3
$ws.Range("K113").Value = 1370.3
$ws.Range("M113").Value = 799.7

$ws.Range("H122").Value = 2733.4
$ws.Range("J122").Value = 3699
$ws.Range("L122").Value = 11097
$ws.Range("N122").Value = -15997

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2756.25
$ws.Range("I122").Value = 2721.5
$ws.Range("K122").Value = 8164.5
$ws.Range("M122").Value = -5714.5

$ws.Range("H131").Value = 137521.42
$ws.Range("J131").Value = 147000
$ws.Range("L131").Value = 147000
$ws.Range("N131").Value = -157080

$ws.Range("H132").Value = 9078.074000000001
$ws.Range("I132").Value = 8497.450000000001
$ws.Range("J132").Value = 10737
$ws.Range("K132").Value = 25492.35
$ws.Range("L132").Value = 32211
$ws.Range("M132").Value = -22962.35
$ws.Range("N132").Value = -37271

$ws.Range("H136").Value = 3544.8215
$ws.Range("I136").Value = 1479.1428
$ws.Range("J136").Value = 5610.5
$ws.Range("K136").Value = 4437.428400000001
$ws.Range("L136").Value = 16831.5
$ws.Range("M136").Value = -1887.428400000001
$ws.Range("N136").Value = -21931.5
